# Add season record columns (Wins / Losses / Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy formatting from the last existing header cell (AB1), then set values
$ws.Range("AB1").Copy($ws.Range("AC1:AE1"))
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Data rows 2-42: same season record repeated for every player row
for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 29).Value = 97
    $ws.Cells.Item($row, 30).Value = 65
    $ws.Cells.Item($row, 31).Value = 0
}
